$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @{
    3  = "28-07-2022"
    4  = "01-08-2022"
    5  = "04-08-2022"
    6  = "08-08-2022"
    7  = "11-08-2022"
    8  = "15-08-2022"
    9  = "18-08-2022"
    10 = "22-08-2022"
    11 = "25-08-2022"
    12 = "29-08-2022"
    13 = "01-09-2022"
    14 = "05-09-2022"
    15 = "08-09-2022"
    16 = "12-09-2022"
    17 = "15-09-2022"
    18 = "19-09-2022"
    19 = "22-09-2022"
    20 = "26-09-2022"
    21 = "29-09-2022"
}

foreach ($row in $dates.Keys) {
    $cell = $ws.Range("A$row")
    # Keep these as plain text (matches inlineStr/text cells in the original
    # file) instead of letting Excel auto-convert the dd-mm-yyyy pattern
    # into a date serial number.
    $cell.NumberFormat = "@"
    $cell.Value = $dates[$row]
}

# Row 13 additional value updates (attendance counted as present instead of absent)
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("H13").Value = 0
